$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Rows 1-3: change the single summary values to "0M"
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"

# Insert 10 new rows right after row 3 (before the current row 4),
# in document order: 147, 0.00003, 0.00026, 0.00008, 0.00005, 0.00007,
# 0.00007, 0.00026, 0.01089, 100.0
$newValues = @('147','0.00003','0.00026','0.00008','0.00005','0.00007','0.00007','0.00026','0.01089','100.0')
$refRow = $t.Rows.Item(4)
$reversedValues = $newValues[($newValues.Length - 1)..0]
foreach ($v in $reversedValues) {
    $newRow = $t.Rows.Add($refRow)
    $newRow.Cells.Item(1).Range.Text = $v
}

# The last three rows of the table previously held tab-separated multi-value
# runs; collapse each back down to a single value.
$rc = $t.Rows.Count
$t.Cell($rc - 2, 1).Range.Text = "100"
$t.Cell($rc - 1, 1).Range.Text = "0.01"
$t.Cell($rc, 1).Range.Text = "228"
